$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row for columns D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), O (Origen), P (Precio $/Kg)
# derived from permuting the original rows (2-14) per the commit diff.

$data = @{
    2  = @{ D = 44232; J = 250; K = 5000; L = 6000; M = 5500; O = "Provincia de Quillota"; P = 344 }
    3  = @{ D = 44230; J = 250; K = 5000; L = 6000; M = 5500; O = "Provincia de Quillota"; P = 344 }
    4  = @{ D = 44187; J = 160; K = 5000; L = 6000; M = 5500; O = "Provincia de Quillota"; P = 344 }
    5  = @{ D = 44231; J = 250; K = 5000; L = 6000; M = 5500; O = "Provincia de Quillota"; P = 344 }
    6  = @{ D = 44251; J = 120; K = 5000; L = 5000; M = 5000; O = "Región Metropolitana"; P = 312 }
    7  = @{ D = 44210; J = 340; K = 5000; L = 6000; M = 5500; O = "Provincia de Quillota"; P = 344 }
    8  = @{ D = 44188; J = 210; K = 5000; L = 6000; M = 5500; O = "Provincia de Quillota"; P = 344 }
    9  = @{ D = 44189; J = 250; K = 5000; L = 6000; M = 5500; O = "Provincia de Quillota"; P = 344 }
    10 = @{ D = 44208; J = 160; K = 5000; L = 6000; M = 5500; O = "Provincia de Quillota"; P = 344 }
    11 = @{ D = 44204; J = 430; K = 5000; L = 6000; M = 5500; O = "Provincia de Quillota"; P = 344 }
    12 = @{ D = 44292; J = 90;  K = 6000; L = 6000; M = 6000; O = "Región Metropolitana"; P = 375 }
    13 = @{ D = 44215; J = 250; K = 5000; L = 6000; M = 5500; O = "Provincia de Quillota"; P = 344 }
    14 = @{ D = 44186; J = 160; K = 5000; L = 6000; M = 5500; O = "Provincia de Quillota"; P = 344 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value  = $vals.D   # D - Fecha
    $ws.Cells.Item($row, 10).Value = $vals.J   # J - Volumen
    $ws.Cells.Item($row, 11).Value = $vals.K   # K - Precio minimo
    $ws.Cells.Item($row, 12).Value = $vals.L   # L - Precio maximo
    $ws.Cells.Item($row, 13).Value = $vals.M   # M - Precio promedio ponderado
    $ws.Cells.Item($row, 15).Value = $vals.O   # O - Origen
    $ws.Cells.Item($row, 16).Value = $vals.P   # P - Precio $/Kg
}
